$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bob took the "Architecture design" task - set C4 (assigned person) to "ehab"
$ws.Range("C4").Value = "ehab"

# Update the active cell selection to G14
$ws.Range("G14").Select()
